$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "as of" date in the confidential disclosure text (A10)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-28 for illustrative purposes only and are subject to change."

# Update weight (D) and percent-change (E) figures for rows 2-7
$ws.Range("D2").Value = 0.2534548305943569
$ws.Range("E2").Value = -0.002251527822450861

$ws.Range("D3").Value = 0.4875079094003634
$ws.Range("E3").Value = -0.002406417112299386

$ws.Range("D4").Value = 0.101579591103425
$ws.Range("E4").Value = -0.006498329001114023

$ws.Range("D5").Value = 0.1000205180945426
$ws.Range("E5").Value = -0.003021148036253751

$ws.Range("D6").Value = 0.05743715080731197
$ws.Range("E6").Value = 0.001992031872509958

$ws.Range("E7").Value = -0.002591665737921023
